# Apply the "output generated at 456a3b4" refresh to the 合肥-漫展信息 workbook.
#
# For the "展览" and "全部类型" sheets the scraped event list was refreshed:
#   - two events that are no longer listed were removed entirely:
#       "合肥·ACGN夏日游园会第七回-泳池派对"            (id=91677)
#       "肥西·星域动漫游戏嘉年华"                        (id=90489)
#   - the remaining events kept their place (shifting up to fill the gap)
#     but their "想去人数" (interest count, column F) was refreshed with an
#     updated number pulled from the source site.

$wb = $excel.ActiveWorkbook

function Update-FanZhanSheet($SheetName, $FValues) {
    $ws = $wb.Worksheets.Item($SheetName)

    # Delete the two rows that disappeared from the listing. Row 4 (肥西·星域动漫游戏嘉年华)
    # is removed first so that row 2 (合肥·ACGN夏日游园会第七回-泳池派对) still refers to the
    # correct row when it is deleted afterwards.
    $ws.Rows.Item(4).Delete()
    $ws.Rows.Item(2).Delete()

    # Refresh the "想去人数" (interest count) column for every remaining data row.
    foreach ($row in ($FValues.Keys | Sort-Object)) {
        $ws.Cells.Item([int]$row, 6).Value2 = $FValues[$row]
    }
}

$sheet1FValues = @{
    2 = 243
    3 = 67
    4 = 13
    5 = 5810
    6 = 5142
    7 = 131
    8 = 57
    9 = 5
    10 = 60
    11 = 217
    12 = 28
}

$sheet4FValues = @{
    2 = 243
    3 = 67
    4 = 13
    5 = 5810
    6 = 5142
    7 = 131
    8 = 57
    9 = 5
    10 = 60
    11 = 217
    12 = 79
    13 = 5
    14 = 28
    15 = 2
}

Update-FanZhanSheet "展览" $sheet1FValues
Update-FanZhanSheet "全部类型" $sheet4FValues
